$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original row values (rows 2,3,4) for the columns that rotate: D, J, K, L, M, O, P
$origD2 = $ws.Range("D2").Value()
$origJ2 = $ws.Range("J2").Value()
$origK2 = $ws.Range("K2").Value()
$origL2 = $ws.Range("L2").Value()
$origM2 = $ws.Range("M2").Value()
$origO2 = $ws.Range("O2").Value()
$origP2 = $ws.Range("P2").Value()

$origD3 = $ws.Range("D3").Value()
$origJ3 = $ws.Range("J3").Value()
$origK3 = $ws.Range("K3").Value()
$origL3 = $ws.Range("L3").Value()
$origM3 = $ws.Range("M3").Value()
$origO3 = $ws.Range("O3").Value()
$origP3 = $ws.Range("P3").Value()

$origD4 = $ws.Range("D4").Value()
$origJ4 = $ws.Range("J4").Value()
$origK4 = $ws.Range("K4").Value()
$origL4 = $ws.Range("L4").Value()
$origM4 = $ws.Range("M4").Value()
$origO4 = $ws.Range("O4").Value()
$origP4 = $ws.Range("P4").Value()

# Row 2 <- original Row 3
$ws.Range("D2").Value = $origD3
$ws.Range("J2").Value = $origJ3
$ws.Range("K2").Value = $origK3
$ws.Range("L2").Value = $origL3
$ws.Range("M2").Value = $origM3
$ws.Range("O2").Value = $origO3
$ws.Range("P2").Value = $origP3

# Row 3 <- original Row 4
$ws.Range("D3").Value = $origD4
$ws.Range("J3").Value = $origJ4
$ws.Range("K3").Value = $origK4
$ws.Range("L3").Value = $origL4
$ws.Range("M3").Value = $origM4
$ws.Range("O3").Value = $origO4
$ws.Range("P3").Value = $origP4

# Row 4 <- original Row 2
$ws.Range("D4").Value = $origD2
$ws.Range("J4").Value = $origJ2
$ws.Range("K4").Value = $origK2
$ws.Range("L4").Value = $origL2
$ws.Range("M4").Value = $origM2
$ws.Range("O4").Value = $origO2
$ws.Range("P4").Value = $origP2
